$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

$updates = @{
    2  = 3000
    3  = 2000
    4  = 0
    5  = 700
    9  = 2500
    10 = 0
    11 = 2500
    12 = 400
    13 = 500
    14 = 2000
    15 = 0
    16 = 1000
    17 = 0
    18 = 2000
    21 = 3500
    24 = 700
    25 = 2500
    26 = 0
    27 = 4000
    28 = 170
    30 = 2000
    31 = 2500
    33 = 1000
    34 = 3000
    37 = 2500
    38 = 40470
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
